$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 83 and 84 (match ids 81/82) had their data swapped between each
# other in the source update - id/Div/Div Original Name/Date (columns
# A, C, D, E) stay put, but everything else (B, F:AC) trades places.

$b83 = $ws.Range("B83").Value2
$b84 = $ws.Range("B84").Value2
$rest83 = $ws.Range("F83:AC83").Value2
$rest84 = $ws.Range("F84:AC84").Value2

$ws.Range("B83").Value2 = $b84
$ws.Range("B84").Value2 = $b83

$ws.Range("F83:AC83").Value2 = $rest84
$ws.Range("F84:AC84").Value2 = $rest83
